# Append 10 new conversation rows (rows 74-83) to Sheet1, matching the
# existing "conversations_with_media" export layout:
#   A: Timestamp (text, "yyyy-MM-dd HH:mm:ss")
#   B: Sender (text)
#   C: Sender Id (number)
#   D: Phone (text - numeric-looking value kept as text)
#   E: Message (text)
#   F: Media (blank)
#   G: Channel (blank)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 74; Timestamp = "2025-10-18 09:15:30"; Sender = "Noah"; SenderId = 8450689526; Phone = "13052054965"; Message = "This is a test message" },
    @{ Row = 75; Timestamp = "2025-10-18 09:16:44"; Sender = "Noah"; SenderId = 8450689526; Phone = "13052054965"; Message = "This is a test message" },
    @{ Row = 76; Timestamp = "2025-10-18 09:19:57"; Sender = "Noah"; SenderId = 8450689526; Phone = "13052054965"; Message = "This is a test message" },
    @{ Row = 77; Timestamp = "2025-10-18 09:23:22"; Sender = "Noah"; SenderId = 8450689526; Phone = "13052054965"; Message = "This is a test message" },
    @{ Row = 78; Timestamp = "2025-10-18 09:25:43"; Sender = "Noah"; SenderId = 8450689526; Phone = "13052054965"; Message = "Test" },
    @{ Row = 79; Timestamp = "2025-10-18 09:40:59"; Sender = "Noah"; SenderId = 8450689526; Phone = "13052054965"; Message = "Test" },
    @{ Row = 80; Timestamp = "2025-10-18 09:43:28"; Sender = "Noah"; SenderId = 8450689526; Phone = "13052054965"; Message = "This message was sent at 9:43 am" },
    @{ Row = 81; Timestamp = "2025-10-18 09:45:59"; Sender = "Noah"; SenderId = 8450689526; Phone = "13052054965"; Message = "This message was sent at 9:45 am" },
    @{ Row = 82; Timestamp = "2025-10-18 10:15:53"; Sender = "Noah"; SenderId = 8450689526; Phone = "13052054965"; Message = "This message was sent at 10:15 am" },
    @{ Row = 83; Timestamp = "2025-10-18 10:18:15"; Sender = "Noah"; SenderId = 8450689526; Phone = "13052054965"; Message = "This message was sent at 10:18 am" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Timestamp
    $ws.Cells.Item($r.Row, 2).Value = $r.Sender
    $ws.Cells.Item($r.Row, 3).Value = $r.SenderId
    # Leading apostrophe forces the numeric-looking phone number to be
    # stored as text (quote-prefix), matching the rest of the sheet.
    $ws.Cells.Item($r.Row, 4).Value = "'" + $r.Phone
    $ws.Cells.Item($r.Row, 5).Value = $r.Message
    # A lone leading apostrophe marks these as empty *text* cells (rather
    # than leaving them as untouched/blank cells) so Media/Channel read
    # back as empty strings, matching the source export.
    $ws.Cells.Item($r.Row, 6).Value = "'"
    $ws.Cells.Item($r.Row, 7).Value = "'"
}
